# Update "想去人数" (want-to-go headcount) values in column F across the
# four worksheets, per the commit "Update gh-pages to output generated at 456a3b4".
#
# Sheet 1: 展览   (Exhibition)
# Sheet 2: 演出   (Performance)
# Sheet 3: 本地生活 (Local Life)
# Sheet 4: 全部类型 (All Types) - combined view of the other three sheets

$wb = $excel.ActiveWorkbook

function Set-FValues($Worksheet, $Updates) {
    foreach ($row in $Updates.Keys) {
        $Worksheet.Cells.Item([int]$row, 6).Value = $Updates[$row]
    }
}

# 展览 (sheet 1)
$ws1 = $wb.Worksheets.Item("展览")
Set-FValues $ws1 @{
    3  = 265
    4  = 1070
    5  = 2413
    7  = 675
    8  = 49
    9  = 223
    10 = 179
    11 = 679
    12 = 78
    13 = 103
    14 = 1432
    15 = 103
    17 = 193
}

# 演出 (sheet 2)
$ws2 = $wb.Worksheets.Item("演出")
Set-FValues $ws2 @{
    10 = 14
    19 = 46
}

# 本地生活 (sheet 3)
$ws3 = $wb.Worksheets.Item("本地生活")
Set-FValues $ws3 @{
    2 = 6338
    3 = 793
    4 = 2000
    5 = 228
}

# 全部类型 (sheet 4) - combined sheet, includes rows mirroring sheets 1-3
$ws4 = $wb.Worksheets.Item("全部类型")
Set-FValues $ws4 @{
    2  = 6338
    3  = 793
    4  = 2000
    5  = 228
    11 = 265
    12 = 1070
    16 = 2413
    19 = 14
    22 = 675
    23 = 49
    24 = 223
    26 = 179
    27 = 679
    28 = 78
    29 = 103
    31 = 1432
    32 = 103
    36 = 193
    39 = 46
}
